$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 21: helloworld6 / Visitor / isValidated = false
$ws.Range("B21").Value = -6
$ws.Range("C21").Value = "helloworld6@gmail.com"
$ws.Range("D21").Value = "helloworld6"
$ws.Range("E21").Value = "Visitor"
$ws.Range("F21").Value = $false
$ws.Range("G21").Value = "'"
$ws.Range("G21").Style = "Normal"
$ws.Range("H21").Value = "'"
$ws.Range("H21").Style = "Normal"
$ws.Range("I21").Value = 0

# Row 22: helloworld7 / Visitor / isValidated = true
$ws.Range("B22").Value = -7
$ws.Range("C22").Value = "helloworld7@gmail.com"
$ws.Range("D22").Value = "helloworld7"
$ws.Range("E22").Value = "Visitor"
$ws.Range("F22").Value = $true
$ws.Range("G22").Value = "'"
$ws.Range("G22").Style = "Normal"
$ws.Range("H22").Value = "'"
$ws.Range("H22").Style = "Normal"
$ws.Range("I22").Value = 0
